# Apply edits described by the diff to Nomina.xlsx

$wb = $excel.ActiveWorkbook
$wsNomina = $wb.Worksheets.Item("Nomina")
$wsCalculo = $wb.Worksheets.Item("Calculo_nomina")

# --- Sheet1 "Nomina" cell content changes ---

# Row 5: H5 message becomes path-specific for periodType
$wsNomina.Range("H5").Value = "is an invalid start of a value. Path: `$.periodType"

# Row 6: H6 message becomes path-specific for year; add D6 = 1
$wsNomina.Range("H6").Value = "is an invalid start of a value. Path: `$.year"
$wsNomina.Range("D6").Value = 1

# Row 7: H7 message becomes path-specific for month; add D7 = 1, F7 = 2021
$wsNomina.Range("H7").Value = "is an invalid start of a value. Path: `$.month"
$wsNomina.Range("D7").Value = 1
$wsNomina.Range("F7").Value = 2021
$wsNomina.Range("P7").Value = 400
$wsNomina.Range("R7").Value = 400

# Row 8: P8/R8 change from 404 to 401
$wsNomina.Range("P8").Value = 401
$wsNomina.Range("R8").Value = 401

# --- Sheet1 view changes ---
$wsNomina.Application.ActiveWindow.ScrollColumn = 14
$wsNomina.Range("R9").Select()

# --- Sheet2 "Calculo_nomina" view changes ---
$wsCalculo.Activate()
$wsCalculo.Range("B39").Select()

# Re-activate Nomina sheet and restore tab selection
$wsNomina.Activate()

# --- Workbook window size ---
$excel.ActiveWindow.Width = 7476
$excel.ActiveWindow.Height = 2232
